# LE module status in settings
# - Add a new row below "Astrophotography countdown fix" showing status
#   of "Max exposure toast" (all columns marked "OK").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 45: A45 = label, B45:F45 = "OK"
$ws.Range("A45").Value = "Max exposure toast"
$ws.Range("B45:F45").Value = "OK"

# Copy the style of the row above (row 44) onto the new row so that
# formatting (blue "note" font) carries over.
$ws.Range("A44:F44").Copy() | Out-Null
$ws.Range("A45:F45").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Move the active selection, matching the workbook's saved cursor position.
$ws.Range("H40").Select() | Out-Null

# Stop printing cell comments at sheet end (now "none").
$ws.PageSetup.PrintComments = "NoComments"

# Slightly widen the sheet's default column width (10.07 -> 10.08 chars).
$ws.StandardWidth = 10.08203125
